$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# NumberFormat "@" (Text) is applied to Price cells before assignment so
# numeric-looking strings (e.g. "322.06", "5.980") are stored verbatim as
# text instead of being parsed into floating point numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.954.73'
$ws.Range("E2").Value = '  -2.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.768.43'
$ws.Range("E3").Value = '  -3.53%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.06'
$ws.Range("E5").Value = '  -2.50%  '
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4267'
$ws.Range("E7").Value = '  -6.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3636'
$ws.Range("E8").Value = '  -4.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.47'
$ws.Range("E9").Value = '  -3.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07461'
$ws.Range("E10").Value = '  -5.00%  '
$ws.Range("E11").Value = '  -4.49%  '
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.09'
$ws.Range("E13").Value = '  -6.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.115'
$ws.Range("E14").Value = '  -4.74%  '
$ws.Range("E15").Value = '  -2.92%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.799.63'
$ws.Range("E16").Value = '  -2.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.18'
$ws.Range("E17").Value = '  -2.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001061'
$ws.Range("E18").Value = '  -2.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06409'
$ws.Range("E19").Value = '  -0.16%  '
$ws.Range("E20").Value = '  +0.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.18'
$ws.Range("E21").Value = '  -2.92%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.980'
$ws.Range("E22").Value = '  -6.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.970.31'
$ws.Range("E23").Value = '  -2.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.37'
$ws.Range("E24").Value = '  -3.73%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.102'
$ws.Range("E25").Value = '  -8.40%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.09'
$ws.Range("E26").Value = '  +2.65%  '
$ws.Range("E27").Value = '  -3.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.996.16'
$ws.Range("E28").Value = '  -2.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.173'
$ws.Range("E29").Value = '  -9.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.78'
$ws.Range("E30").Value = '  -3.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.153'
$ws.Range("E31").Value = '  -5.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.740'
$ws.Range("E32").Value = '  +1.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.629'
$ws.Range("E33").Value = '  -5.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08895'
$ws.Range("E34").Value = '  -4.95%  '
$ws.Range("E35").Value = '  -4.81%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02311'
$ws.Range("E36").Value = '  -2.71%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2120'
$ws.Range("E37").Value = '  -4.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.037'
$ws.Range("E38").Value = '  -4.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06027'
$ws.Range("E39").Value = '  -4.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6372'
$ws.Range("E40").Value = '  -5.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.185'
$ws.Range("E41").Value = '  -1.20%  '
$ws.Range("E42").Value = '  +0.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.403'
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.842'
$ws.Range("E44").Value = '  -4.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.45'
$ws.Range("E45").Value = '  -5.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5941'
$ws.Range("E46").Value = '  -4.18%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.699'
$ws.Range("E47").Value = '  -2.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.007'
$ws.Range("E48").Value = '  -3.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '122.74'
$ws.Range("E49").Value = '  -4.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.189'
$ws.Range("E50").Value = '  +2.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06860'
$ws.Range("E51").Value = '  -2.43%  '
